$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing imae (column B) values for rows 218-226
$ws.Cells.Item(218, 2).Value = 125.595860670832
$ws.Cells.Item(219, 2).Value = 130.344813000994
$ws.Cells.Item(220, 2).Value = 126.12797294583
$ws.Cells.Item(221, 2).Value = 122.845678576047
$ws.Cells.Item(222, 2).Value = 132.913461098593
$ws.Cells.Item(223, 2).Value = 124.801835203815
$ws.Cells.Item(224, 2).Value = 126.289707999204
$ws.Cells.Item(225, 2).Value = 125.764151973952
$ws.Cells.Item(226, 2).Value = 122.484079339617

# Add a new row 227 (periodo 2025-10-01) with its data, copying the
# date formatting used by the rest of column A
$ws.Range("A226").Copy()
$ws.Range("A227").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A227").Value = 45931
$ws.Range("B227").Value = 123.176380654678
$ws.Range("C227").Value = 172.26
$ws.Range("D227").Value = 149.92

$excel.CutCopyMode = 0
